$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "43.072.44"
$ws.Range("E2").Value = "  +0.04%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.305.37"
$ws.Range("E3").Value = "  -0.31%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 (BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "

# Row 6 (Solana)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.57%  "

# Row 7 (XRP)
$ws.Range("E7").Value = "  +3.79%  "

# Row 8 (USDC)
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 (Cardano)
$ws.Range("E9").Value = "  +1.19%  "

# Row 10 (Avalanche)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.36%  "

# Row 11 (Dogecoin)
$ws.Range("E11").Value = "  -0.14%  "

# Row 12 (TRON)
$ws.Range("E12").Value = "  -0.40%  "

# Row 13 (Chainlink)
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.05%  "

# Row 14 (Polkadot)
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.46%  "

# Row 15 (WrappedliquidstakedEther2.0)
$ws.Range("D15").Value = "2.663.48"
$ws.Range("E15").Value = "  -0.35%  "

# Row 16 (WrappedEther)
$ws.Range("D16").Value = "2.350.04"
$ws.Range("E16").Value = "  +1.77%  "

# Row 17 (Polygon)
$ws.Range("E17").Value = "  -2.23%  "

# Row 18 (WrappedBTC)
$ws.Range("D18").Value = "42.974.99"
$ws.Range("E18").Value = "  +0.06%  "

# Row 19 (InternetComputer(DFINITY))
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.88%  "

# Row 20 (ShibaInu)
$ws.Range("E20").Value = "  +0.63%  "

# Row 21 (Uniswap)
$ws.Range("E21").Value = "  -1.28%  "

# Row 22 (Litecoin)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.71%  "

# Row 23 (BitcoinCash)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.20%  "

# Row 24 (ImmutableX)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.81%  "

# Row 25 (Dai)
$ws.Range("E25").Value = "  -0.07%  "

# Row 26 (PancakeSwap)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.00%  "

# Row 27 (EthereumClassic)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.27%  "

# Row 28 (Monero)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.71%  "

# Row 29 (Cosmos)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "

# Row 30 (Toncoin)
$ws.Range("E30").Value = "  -13.28%  "

# Row 31 (InjectiveProtocol)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.96%  "

# Row 32 (Filecoin)
$ws.Range("E32").Value = "  +4.31%  "

# Row 33 (FirstDigitalUSD)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.03%  "

# Row 34 (RenderToken)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.63%  "

# Row 35 (Celestia)
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.19%  "

# Row 36 (WEMIXToken)
$ws.Range("E36").Value = "  -0.55%  "

# Row 37 (Hedera)
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0691"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.44%  "

# Row 38 (Kaspa)
$ws.Range("E38").Value = "  -0.65%  "

# Row 39 (ARBITRUM)
$ws.Range("E39").Value = "  +0.32%  "

# Row 40 (Stellar)
$ws.Range("E40").Value = "  +2.09%  "

# Row 41 (LidoDAOToken)
$ws.Range("E41").Value = "  -3.20%  "

# Row 42 (Maker)
$ws.Range("D42").Value = "2.004.29"
$ws.Range("E42").Value = "  +0.91%  "

# Row 43 (VeChain)
$ws.Range("E43").Value = "  -0.23%  "

# Row 44 (ApeXProtocol)
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.53%  "

# Row 45 (FraxShare)
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.84%  "

# Row 46 (EnergySwap)
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.66%  "

# Row 47 (NEARProtocol)
$ws.Range("E47").Value = "  -2.99%  "

# Row 48 (MultiversX)
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.82%  "

# Row 49 (RocketPoolETH)
$ws.Range("D49").Value = "2.529.56"
$ws.Range("E49").Value = "  -0.42%  "

# Row 50 (BitcoinSV)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.25%  "

# Row 51 (Stacks)
$ws.Range("E51").Value = "  +0.46%  "
